# retire ":" dans freq
# Update the "Fréquences" (D) and "Taux total" (F) columns with new
# computed values, drop the now-stale hidden "_xlchart.v1.*" defined
# names left over from a prior chart edit, and restore the sheet's last
# selection. The remaining cosmetic artifacts (window size/position,
# chart numCache mirrors, per-point ext blocks, dLblPos, revision ids,
# …) are internal Excel bookkeeping that is re-derived automatically
# from the sheet data and isn't something a user-level edit sets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Drop the stale hidden "_xlchart.v1.*" helper names Excel had generated
# for the chart's "filtered series" cache — they no longer referenced the
# current chart state and were cleaned up.
$namesToDelete = @()
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $namesToDelete += $wb.Names.Item($i).Name
}
foreach ($n in $namesToDelete) {
    $wb.Names.Item($n).Delete()
}

$ws.Range("D2").Value = 36
$ws.Range("F2").Value = -3.4443999999999999

$ws.Range("D3").Value = 21
$ws.Range("F3").Value = -4.5

$ws.Range("D4").Value = 82
$ws.Range("F4").Value = -1.75

$ws.Range("D5").Value = 406
$ws.Range("F5").Value = 0.41539999999999999

$ws.Range("D6").Value = 507
$ws.Range("F6").Value = 0.43919999999999998

# Update the active selection to match the post-edit state recorded in
# the workbook (B2:F6 selected, active cell F6).
$ws.Range("B2:F6").Select() | Out-Null
